$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.290.32'
$ws.Range("E2").Value = '  -0.11%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.867.79'
$ws.Range("E3").Value = '  +0.35%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7246'
$ws.Range("E5").Value = '  +3.16%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.98'
$ws.Range("E6").Value = '  +1.22%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07876'
$ws.Range("E8").Value = '  +0.49%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3093'
$ws.Range("E9").Value = '  +1.46%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.23'
$ws.Range("E10").Value = '  +2.05%  '

# Row 11
$ws.Range("E11").Value = '  +1.20%  '

# Row 12
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7229'
$ws.Range("E12").Value = '  +1.34%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.241'
$ws.Range("E13").Value = '  +0.43%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.849.11'
$ws.Range("E14").Value = '  -1.13%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.75'
$ws.Range("E15").Value = '  +1.81%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.346.56'
$ws.Range("E16").Value = '  -0.14%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.864'
$ws.Range("E17").Value = '  +1.05%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '243.81'
$ws.Range("E18").Value = '  +1.99%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007834'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.23'
$ws.Range("E20").Value = '  +0.37%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.120.82'
$ws.Range("E21").Value = '  -0.22%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9997'
$ws.Range("E22").Value = '  -0.14%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.024'
$ws.Range("E23").Value = '  +6.84%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  +0.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1599'
$ws.Range("E25").Value = '  +12.49%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.38'
$ws.Range("E26").Value = '  -0.18%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.966'
$ws.Range("E27").Value = '  +0.77%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.26'
$ws.Range("E28").Value = '  +1.08%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.344'
$ws.Range("E29").Value = '  -2.32%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.495'
$ws.Range("E30").Value = '  +1.63%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.400'
$ws.Range("E31").Value = '  +2.47%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.114'
$ws.Range("E32").Value = '  +1.94%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05214'
$ws.Range("E33").Value = '  +0.85%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.935'
$ws.Range("E34").Value = '  +1.71%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.186'
$ws.Range("E35").Value = '  +0.67%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7288'
$ws.Range("E36").Value = '  +3.37%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.681'
$ws.Range("E37").Value = '  +0.14%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01857'
$ws.Range("E38").Value = '  +0.71%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.702'
$ws.Range("E39").Value = '  +0.29%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.167.95'
$ws.Range("E40").Value = '  -0.46%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9053'
$ws.Range("E41").Value = '  -1.42%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.101'
$ws.Range("E42").Value = '  +1.37%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.30'
$ws.Range("E43").Value = '  +1.06%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.06%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.87'
$ws.Range("E45").Value = '  +0.05%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.016.42'
$ws.Range("E46").Value = '  -0.24%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5282'
$ws.Range("E47").Value = '  -1.27%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.776'
$ws.Range("E48").Value = '  +1.38%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.895'
$ws.Range("E49").Value = '  +6.29%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.277'
$ws.Range("E50").Value = '  +1.41%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4273'
$ws.Range("E51").Value = '  +0.59%  '
